$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 10 (weekly price update: two new
# "Ciboulette" quotes for 2023-10-?? get added at the top of the historical
# list below the most recent entries, pushing the rest of the series down).
$ws.Rows("10:11").Insert()

# Row 10: new "Primera" quote
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 45203
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 100112039
$ws.Range("G10").Value = "Ciboulette"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 2500
$ws.Range("L10").Value = 2500
$ws.Range("M10").Value = 2500
$ws.Range("N10").Value = "$/docena de atados"
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value = 833
$ws.Range("Q10").Value = 3
$ws.Range("R10").Value = "Hortaliza"

# Row 11: new "Segunda" quote
$ws.Range("A11").Value = 7
$ws.Range("B11").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C11").Value = "Ñuble"
$ws.Range("D11").Value = 45203
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = 100112039
$ws.Range("G11").Value = "Ciboulette"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 1500
$ws.Range("L11").Value = 1500
$ws.Range("M11").Value = 1500
$ws.Range("N11").Value = "$/docena de atados"
$ws.Range("O11").Value = "Región Metropolitana"
$ws.Range("P11").Value = 500
$ws.Range("Q11").Value = 3
$ws.Range("R11").Value = "Hortaliza"
